# [CI] Auto-release exports files
# Refresh the generated KiCost report's timestamp/version stamps and the
# recalculated EUR/USD rate, matching the upstream KiCost v1.1.10 re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Sun Sep  4 19:50:10 2022"   # Prj date:
$ws.Range("B4").Value = "2022-09-04 19:50:12"        # $ date:
$ws.Range("A13").Value = "KiCost® v1.1.10"           # KiCost version footer
$ws.Range("C11").Value = 1.00070049034324            # EUR(€)/USD($) rate
